$d = $word.ActiveDocument

# 1. Font rename: TimesNewToman -> Times New Roman across the whole document.
#    (Find/Replace only touches visible text, not rFonts attributes, so we
#    set Font.Name on a full-document Range instead.)
$fullRange = $d.Range(0, $d.Content.End)
$fullRange.Font.Name = "Times New Roman"

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Text = $old
    $find.Replacement.Text = $new
    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null
}

# Title
Replace-Text "Celestial Harmony: Unraveling the Cosmic Symphony" "History: A Window to the Past, a Path to the Future"

# Author line ("Dr" + "." + " Amelia Rodriguez" -> single run "Earl Watson")
Replace-Text "Dr. Amelia Rodriguez" "Earl Watson"

# Email user part
Replace-Text "ameliarodriguez@cosmosobservatory" "newwatson@hope-mail"

# Body paragraph sentences
Replace-Text "Lost in the depths of the cosmos lies an enigmatic symphony, a celestial harmony that has captivated the imaginations of astronomers, physicists, and philosophers for millennia" "History is the tapestry of human existence, a chronicle of our triumphs and follies, our hopes and dreams"

Replace-Text " From the gentle hum of primordial atoms to the thunderous roar of exploding stars, the universe resonates with a mesmerizing symphony that holds the key to unlocking some of the universe's most profound mysteries" " It's an intricate web of stories, each thread telling a unique tale of ambition, resilience, and the never-ending quest for progress"

Replace-Text " In this essay, we will delve into the fascinating realm of cosmic music, exploring the intricacies of this celestial symphony and its profound implications for our understanding of the universe" " As students, we have the privilege of exploring this vast tapestry, uncovering the hidden patterns and gaining a deeper understanding of ourselves and our place in the universe"

Replace-Text "The harmonious fusion of diverse celestial objects, each contributing its unique melody, composes the cosmic symphony" "Our ancestors whispered their tales to the winds, inscribing them on cave walls, papyrus scrolls, and weathered tombstones"

Replace-Text " Stars, galaxies, and cosmic dust dance in intricate patterns, their gravitational interactions creating a symphony of cosmic ballet" " Each generation adds its own brushstroke to the canvas, leaving behind clues to the mysteries of the past"

Replace-Text " As celestial bodies orbit, collide, and explode, they emit a myriad of electromagnetic waves, ranging from gentle radio waves to penetrating gamma rays" " Through history, we witness the evolution of civilization, from the dawn of humanity to the complexities of the modern world"

Replace-Text " These waves, stretching across the vast expanse of the universe, form a cosmic symphony that spans the entire electromagnetic spectrum" " It's a journey of discovery, a voyage across oceans of time, where we encounter kings and queens, heroes and villains, visionaries and dreamers"

Replace-Text "The symphony of the cosmos not only provides aesthetic wonder but also serves as a valuable scientific tool" "History has the power to ignite our imaginations, to transport us to distant lands and bygone eras"

Replace-Text " By deciphering the intricate melodies of the universe, astronomers can glean insights into the properties of celestial objects, their interactions, and the fundamental laws governing the cosmos" " It's a treasure trove of wisdom, a testament to the unyielding spirit of humanity"

Replace-Text " The cosmic symphony reverberates with information, whispering tales of cosmic evolution, black hole singularities, and the fabric of space-time itself" " By studying history, we not only learn about the past but also gain insights into the present, helping us make informed decisions and navigate the challenges of the future"

# Summary paragraph sentences
Replace-Text "The cosmic symphony, an intricate fusion of celestial melodies, offers a window into the enigmatic workings of the universe" "In this journey through history, we've explored the significance of studying the past, revealing the profound impact it has on shaping our identities and understanding of the world"

Replace-Text " From the primordial hum of the universe's birth to the thunderous crescendo of stellar explosions, the cosmos resonates with a symphony that holds the secrets to understanding the universe's evolution and fundamental laws" " History invites us to contemplate the interconnectedness of all things, showcasing the intricate web that weaves together different cultures, periods, and individuals"

Replace-Text " The harmonious interplay of celestial objects, each contributing its unique melody to the cosmic orchestra, provides a mesmerizing spectacle that continues to captivate and inspire scientific inquiry, revealing the profound interconnectedness of all things in the vast expanse of the cosmos" " Through history, we discover the essence of our humanity, the tapestry of triumphs and tragedies that makes us who we are today"

# Append two new sentences (with their own trailing-period runs) to the end
# of the Summary paragraph, matching the formatting of the paragraph's last run.
$summaryPara = $d.Paragraphs($d.Paragraphs.Count)
$endRange = $summaryPara.Range
$endRange.Collapse(0)
$endRange.MoveEnd(1, -1) | Out-Null
$endRange.InsertAfter(" It's a journey that continues to unfold, a story that's forever in the making, with each new generation adding its own chapters.")

# New empty paragraph at the very end of the document body.
$tailRange = $d.Content
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()
